# Notified_Production_Wind.xlsx update: shift the timestamp column (A) forward
# by 11 days (one full data window -> next window, day 45810 -> day 45821) and
# replace the Notified Production values (B) with the newly retrained values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Notified Production (MW)" values for rows 2..97 (B2:B97), in order.
$newValues = @(
    317, 315, 313, 311, 317, 315, 313, 310, 334, 333,
    333, 333, 369, 369, 370, 371, 388, 389, 390, 391,
    376, 375, 374, 373, 324, 323, 323, 323, 286, 287,
    287, 287, 276, 276, 276, 276, 272, 272, 272, 272,
    264, 265, 266, 266, 297, 298, 299, 299, 357, 358,
    359, 359, 420, 422, 423, 423, 494, 495, 497, 497,
    583, 585, 586, 588, 621, 623, 625, 626, 676, 676,
    676, 675, 709, 708, 707, 707, 758, 759, 761, 762,
    806, 807, 807, 807, 818, 818, 819, 820, 851, 851,
    850, 850, 0, 0, 0, 0
)

$firstRow = 2
$lastRow = 97
$dayShift = 11

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $cellA.Value2 + $dayShift

    $ws.Cells.Item($row, 2).Value = $newValues[$row - $firstRow]
}
